# Adds a "WebSocket Command" actions section to the "Commands" sheet
# (new rows 86-87, mirroring the existing "<X> Command" header row +
# standard sleep()/delay note row pattern used throughout the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

$ws.Range("A86").Value = "WebSocket Command"
$ws.Range("B86").Value = "send(<json>)"
$ws.Range("C86").Value = "If {} substitutions are used, json brackets need to be duplicated to escape them like in send({{ “value”: {}}})"

$ws.Range("B87").Value = "sleep(<float>)"
$ws.Range("C87").Value = "sleep: add a delay of <float> seconds"

[void]$ws.Activate()
[void]$ws.Range("A86").Select()
